# Removal of USDT Tether Omni - Reminder email to clients.docx
# Simplified Chinese (zh) -> Traditional Chinese (zh-Hant) translation update.
#
# The runtime's Word COM-interop layer reliably supports Find/Replace on the
# main document story (body text, including text inside tables and
# hyperlinks). Each replacement below targets a string that appears exactly
# once in the document, so a plain literal Find/Replace (wildcards off) is
# unambiguous and safe.

$d = $word.ActiveDocument

function Replace-Once($find, $replace) {
    $found = $d.Content.Find.Execute(
        $find,     # FindText
        $true,     # MatchCase
        $false,    # MatchWholeWord
        $false,    # MatchWildcards
        $false,    # MatchSoundsLike
        $false,    # MatchAllWordForms
        $true,     # Forward
        1,         # Wrap (wdFindContinue)
        $false,    # Format
        $replace,  # ReplaceWith
        2          # Replace (wdReplaceOne)
    )
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $find)
    } else {
        Write-Output ("OK: " + $find.Substring(0, [Math]::Min(6,$find.Length)))
    }
}

# Title line
Replace-Once "向 ROW 客户发送提醒电子邮件" "向 ROW 客戶傳送提醒電子郵件"

# Subject line
Replace-Once "将于 9 月 29 日移除 Tether Omni (USDT)" "將於 9 月 29 日移除 Tether Omni (USDT)"

# Heading inside body table
Replace-Once "向 Tether Omni 说再见" "向 Tether Omni 道別"

# Main paragraph about GMT cut-off
Replace-Once "自 2023 年 9 月 29 日格林威治标准时间 00:00 起，Deriv 将停止提供Tether Omni (USDT) 作为账户货币。 这是因为 Tether 已停止支持 USDT 的 Omni 转账。" "自 2023 年 9 月 29 日格林威治標準時間 00:00 起，Deriv 將停止提供Tether Omni (USDT) 作為帳戶貨幣。 這是因為 Tether 已停止支援 USDT 的 Omni 轉帳。"

# "What do you need to do?" heading (bold run)
Replace-Once "需要做什么？" "需要做什麼？"

# "If the USDT account " run (immediately followed by the [account ID] run)
Replace-Once "如果 USDT 账户 " "如果 USDT 帳戶 "

# " has a balance, please withdraw ..." run (immediately after [account ID])
Replace-Once " 中有余额，请在上述日期之前提取余额。 如果有持仓头寸，提取余额之前请先平仓。" " 中有餘額，請在上述日期之前提取餘額。 若有持倉頭寸，提取餘額前請先平倉。"

# Hyperlink display text "Check account"
Replace-Once "查看账户" "檢查帳戶"

# Important notice paragraph (account closes / balances transferred)
Replace-Once "USDT 账户将于 2023 年 9 月 29 日格林尼治标准时间 00:00 关闭。 任何持仓头寸将在上述日期后自动平仓，账户余额将转移到最后活跃的账户" "USDT 帳戶將於 2023 年 9 月 29 日格林尼治標準時間 00:00 關閉。 任何持倉頭寸將在上述日期後自動平倉，帳戶餘額將轉移到最後活躍的帳戶"

# Trailing sentence about standard rates/fees (after the comment references)
Replace-Once "在此过程中将采用标准汇率和费用。" "在此過程中將採用標準匯率和費用。"

# "If you have any questions, contact us" paragraph
Replace-Once "如有任何疑问，请通过以下方式联系我们：" "如有任何疑問，請透過以下方式聯繫我們："

# "Live chat" hyperlink text
Replace-Once "实时聊天" "即時聊天"

Write-Output "Done."
